# Append a new "Data in storage:" section to the end of the document,
# right after the existing "Finally, we will make use of tokens..."
# paragraph (the end of the "Data in transfer security" discussion) and
# before the section properties (w:sectPr).
#
# The new content consists of:
#   - one blank spacer paragraph
#   - a bold heading paragraph: "Data in storage:"
#   - a paragraph about encrypting data at rest (contains a mid-run
#     w:lastRenderedPageBreak marker, matching the authored document)
#   - a paragraph about access control
#   - a paragraph about backing up data (contains a proofErr
#     gramStart/gramEnd wrapped "it", matching the authored document)
#
# All of the new runs use the same Times New Roman / en-GB formatting as
# the rest of the document's body text, and the heading run additionally
# carries bold (w:b / w:bCs).
#
# We build the exact WordprocessingML for these paragraphs and drop it in
# with Range.InsertXML so the resulting markup (rPr/pPr ordering, bCs on
# the bold paragraph mark, proofErr tags, the page-break marker, and the
# xml:space="preserve" runs) matches character-for-character rather than
# relying on Find/Replace + ad-hoc formatting calls.

$d = $word.ActiveDocument

$lastParagraph = $d.Paragraphs.Item($d.Paragraphs.Count)
$anchorText = "If the token hasn" + [char]0x2019 + "t expired"
if ($lastParagraph.Range.Text.IndexOf($anchorText) -lt 0) {
    throw "Could not find the expected trailing paragraph to anchor the new content on."
}

$newBodyXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:lang w:val="en-GB"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:bCs/><w:lang w:val="en-GB"/></w:rPr><w:t>Data in storage:</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">While in storage, we will be encrypting data to make it more secure. It’s important that we encrypt personal data such as emails and passwords as we don’t want such information to end up leaked or </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:lang w:val="en-GB"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">stolen. Other than specifically encrypting, we can keep some of the data secure using methods such as password hashing and salting. </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">Another way of keeping the data safe while in storage is by enforcing strict access control. By using access control, we only allow certain authenticated users access to the data and decide which of these users will be </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:lang w:val="en-GB"/></w:rPr><w:t>able to modify this data. We will only allow each user the exact permissions that they will need to carry out their task and nothing more than that. We will be basing the access control on user roles so that each user has the appropriate permissions for their role and responsibilities.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">We will also be backing up the data frequently to ensure that we maintain data integrity. There is always a chance that data can become corrupted, or a problem will arise with </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:lang w:val="en-GB"/></w:rPr><w:t>it</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> and we want to have a solution for this. By backing up the data, we will easily be able to restore it when something like this happens. This will make data recovery easier and will stop us from losing data in these events while also improving the integrity of the data.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$insertionPoint = $d.Content
$insertionPoint.Collapse(0)
[void]$insertionPoint.InsertXML($newBodyXml)
